$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 91.28570999999999
$ws.Range("I5").Value = 91.28570999999999
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 91.28570999999999
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 23.71429000000001

$ws.Range("N8").ClearContents()
$ws.Range("H8").Value = 13.5
$ws.Range("I8").Value = 13.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 40.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 98.5

$ws.Range("H32").Value = 7807.9165
$ws.Range("I32").Value = 6000
$ws.Range("J32").Value = 8169.5
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 8169.5
$ws.Range("M32").Value = -5674
$ws.Range("N32").Value = -8821.5

$ws.Range("H43").Value = 1012.25
$ws.Range("I43").Value = 1075
$ws.Range("J43").Value = 949.5
$ws.Range("K43").Value = 1075
$ws.Range("L43").Value = 949.5
$ws.Range("M43").Value = -1006
$ws.Range("N43").Value = -1087.5

$ws.Range("H48").Value = 1064.6666
$ws.Range("I48").Value = 447.5
$ws.Range("J48").Value = 2299
$ws.Range("K48").Value = 1342.5
$ws.Range("L48").Value = 6897
$ws.Range("M48").Value = -1050.5
$ws.Range("N48").Value = -7481

$ws.Range("H56").Value = 1064.6666
$ws.Range("I56").Value = 447.5
$ws.Range("J56").Value = 2299
$ws.Range("K56").Value = 1342.5
$ws.Range("L56").Value = 6897
$ws.Range("M56").Value = -808.5
$ws.Range("N56").Value = -7965

$ws.Range("H97").Value = 3254.6667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3254.6667
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 9764.000100000001
$ws.Range("N97").Value = -10756.0001

$ws.Range("H98").Value = 696.3333
$ws.Range("I98").Value = 676
$ws.Range("J98").Value = 798
$ws.Range("K98").Value = 676
$ws.Range("L98").Value = 798
$ws.Range("M98").Value = 822
$ws.Range("N98").Value = -3794

$ws.Range("H99").Value = 619.3333
$ws.Range("I99").Value = 619.3333
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1857.9999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -359.9999

$ws.Range("H101").Value = 1724.75
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 1724.75
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 5174.25
$ws.Range("N101").Value = -8418.25

$ws.Range("H122").Value = 696.3333
$ws.Range("I122").Value = 676
$ws.Range("J122").Value = 798
$ws.Range("K122").Value = 2028
$ws.Range("L122").Value = 2394
$ws.Range("M122").Value = 422
$ws.Range("N122").Value = -7294

$ws.Range("H138").Value = 3271.1482
$ws.Range("I138").Value = 1306
$ws.Range("J138").Value = 3958.95
$ws.Range("K138").Value = 3918
$ws.Range("L138").Value = 11876.85
$ws.Range("M138").Value = 1222

$ws.Range("H141").Value = 5749.1055
$ws.Range("I141").Value = 6026.8335
$ws.Range("J141").Value = 750
$ws.Range("K141").Value = 18080.5005
$ws.Range("L141").Value = 2250
$ws.Range("M141").Value = -12900.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3100
$ws.Range("I102").Value = 3100
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3100
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1478

$ws.Range("H122").Value = 6051
$ws.Range("I122").Value = 6387.0454
$ws.Range("J122").Value = 4994.857
$ws.Range("K122").Value = 19161.1362
$ws.Range("L122").Value = 14984.571
$ws.Range("M122").Value = -16711.1362
$ws.Range("N122").Value = -19884.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 132.11111
$ws.Range("I7").Value = 132.11111
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 132.11111
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -19.11111

$ws.Range("H86").Value = 10249.75
$ws.Range("I86").Value = 7999
$ws.Range("J86").Value = 11000
$ws.Range("K86").Value = 7999
$ws.Range("L86").Value = 11000
$ws.Range("M86").Value = -6876
$ws.Range("N86").Value = -13246

$ws.Range("H89").Value = 10249.75
$ws.Range("I89").Value = 7999
$ws.Range("J89").Value = 11000
$ws.Range("K89").Value = 39995
$ws.Range("L89").Value = 55000
$ws.Range("M89").Value = -34379
$ws.Range("N89").Value = -66232

$ws.Range("H122").Value = 957.0769
$ws.Range("I122").Value = 972.6667
$ws.Range("J122").Value = 922
$ws.Range("K122").Value = 2918.0001
$ws.Range("L122").Value = 2766
$ws.Range("M122").Value = -468.0001000000002
$ws.Range("N122").Value = -7666

$ws.Range("H132").Value = 2916
$ws.Range("I132").Value = 1840
$ws.Range("J132").Value = 3633.3333
$ws.Range("K132").Value = 5520
$ws.Range("L132").Value = 10899.9999
$ws.Range("M132").Value = -2990

$ws.Range("H134").Value = 2371.3125
$ws.Range("I134").Value = 2304
$ws.Range("J134").Value = 2573.25
$ws.Range("K134").Value = 6912
$ws.Range("L134").Value = 7719.75
$ws.Range("M134").Value = -4377

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 4110.6665
$ws.Range("I60").Value = 3333
$ws.Range("J60").Value = 4499.5
$ws.Range("K60").Value = 9999
$ws.Range("L60").Value = 13498.5
$ws.Range("M60").Value = -9748
$ws.Range("N60").Value = -14000.5

$ws.Range("H92").Value = 921.125
$ws.Range("I92").Value = 1011.2857
$ws.Range("J92").Value = 290
$ws.Range("K92").Value = 3033.8571
$ws.Range("L92").Value = 870
$ws.Range("M92").Value = -1785.8571

$ws.Range("H94").Value = 5812
$ws.Range("I94").Value = 5812
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 17436
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -16760

$ws.Range("H97").Value = 174.6
$ws.Range("I97").Value = 126.5
$ws.Range("J97").Value = 206.66667
$ws.Range("K97").Value = 379.5
$ws.Range("L97").Value = 620.00001
$ws.Range("M97").Value = 116.5
$ws.Range("N97").Value = -1612.00001

$ws.Range("H103").Value = 3587
$ws.Range("I103").Value = 3105
$ws.Range("J103").Value = 3931.2856
$ws.Range("K103").Value = 9315
$ws.Range("L103").Value = 11793.8568
$ws.Range("M103").Value = -8436
$ws.Range("N103").Value = -13551.8568

$ws.Range("H109").Value = 1630.6
$ws.Range("I109").Value = 1630.6
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 4891.799999999999
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -3851.799999999999

$ws.Range("M131").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4000
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4000
$ws.Range("N73").Value = -5872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6710.222
$ws.Range("I7").Value = 6484.5713
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 6484.5713
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -6372.5713

$ws.Range("H9").Value = 903.5
$ws.Range("I9").Value = 808
$ws.Range("J9").Value = 999
$ws.Range("K9").Value = 808
$ws.Range("L9").Value = 999
$ws.Range("M9").Value = -584

$ws.Range("H40").Value = 6156
$ws.Range("I40").Value = 6156
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6156
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6020

$ws.Range("H126").Value = 6710.222
$ws.Range("I126").Value = 6484.5713
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 19453.7139
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -16983.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M38").ClearContents()
$ws.Range("H38").Value = 23500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 23500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 23500
$ws.Range("N38").Value = -24446

$ws.Range("H52").Value = 3356173.8
$ws.Range("I52").Value = 5014510.5
$ws.Range("J52").Value = 39500
$ws.Range("K52").Value = 5014510.5
$ws.Range("L52").Value = 39500
$ws.Range("M52").Value = -5014284.5
$ws.Range("N52").Value = -39952

$ws.Range("M55").ClearContents()
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0

$ws.Range("H81").Value = 1431
$ws.Range("I81").Value = 984.875
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 1969.75
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -908.75

$ws.Range("H84").Value = 1431
$ws.Range("I84").Value = 984.875
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 9848.75
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -4544.75

$ws.Range("H122").Value = 1977.4
$ws.Range("I122").Value = 1697.1111
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 5091.3333
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -2641.3333
